$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.093.03"

# Row 3
$ws.Range("D3").Value = "2.553.33"
$ws.Range("E3").Value = "  +0.42%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.96"
$ws.Range("E5").Value = "  +2.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.26"
$ws.Range("E6").Value = "  -2.19%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("E9").Value = "  +0.28%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.58"
$ws.Range("E10").Value = "  -1.76%  "

# Row 11
$ws.Range("E11").Value = "  -0.09%  "

# Row 12
$ws.Range("E12").Value = "  -0.95%  "

# Row 13
$ws.Range("E13").Value = "  -2.33%  "

# Row 14
$ws.Range("D14").Value = "3.008.58"
$ws.Range("E14").Value = "  +0.38%  "

# Row 15
$ws.Range("D15").Value = "63.001.84"
$ws.Range("E15").Value = "  -0.51%  "

# Row 16
$ws.Range("E16").Value = "  +1.16%  "

# Row 17
$ws.Range("D17").Value = "2.555.67"
$ws.Range("E17").Value = "  +1.72%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.34"
$ws.Range("E18").Value = "  -2.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.27"
$ws.Range("E19").Value = "  -0.09%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.33"
$ws.Range("E20").Value = "  +0.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.76"
$ws.Range("E21").Value = "  -0.23%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.04%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.56"
$ws.Range("E23").Value = "  -0.52%  "

# Row 24
$ws.Range("D24").Value = "2.676.53"
$ws.Range("E24").Value = "  +0.49%  "

# Row 25
$ws.Range("E25").Value = "  +0.50%  "

# Row 26
$ws.Range("E26").Value = "  +1.21%  "

# Row 27
$ws.Range("B27").Value = "SuiNetwork"
$ws.Range("C27").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.49"
$ws.Range("E27").Value = "  -4.39%  "

# Row 28
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.26%  "

# Row 29
$ws.Range("E29").Value = "  -1.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.69"
$ws.Range("E30").Value = "  +7.67%  "

# Row 31
$ws.Range("E31").Value = "  +5.25%  "

# Row 32
$ws.Range("E32").Value = "  +0.08%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "178.06"
$ws.Range("E33").Value = "  +0.16%  "

# Row 34
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.56"
$ws.Range("E34").Value = "  -0.76%  "

# Row 35
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "421.30"
$ws.Range("E35").Value = "  +0.09%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.401"
$ws.Range("E36").Value = "  -0.76%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.13"
$ws.Range("E37").Value = "  +0.64%  "

# Row 38
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.40"
$ws.Range("E38").Value = "  -0.39%  "

# Row 39
$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.03%  "

# Row 40
$ws.Range("E40").Value = "  -1.26%  "

# Row 41
$ws.Range("E41").Value = "  +0.00%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.71"
$ws.Range("E42").Value = "  +0.61%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "150.79"
$ws.Range("E43").Value = "  -1.84%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.78"
$ws.Range("E44").Value = "  +0.18%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.81"
$ws.Range("E45").Value = "  +0.24%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0538"
$ws.Range("E46").Value = "  +2.91%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.603"
$ws.Range("E47").Value = "  -0.78%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0971"
$ws.Range("E48").Value = "  +0.85%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0239"
$ws.Range("E49").Value = "  +0.53%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.29"
$ws.Range("E50").Value = "  -0.88%  "

# Row 51
$ws.Range("E51").Value = "  -4.68%  "
